# Update Leve profit-tracking figures across all job sheets
# (scheduled market-data refresh for Mandragora_Profits)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 563.7692
$ws.Range("J28").Value = 1253
$ws.Range("L28").Value = 1253
$ws.Range("N28").Value = -2223
# Row 134
$ws.Range("H134").Value = 90948010
$ws.Range("J134").Value = 90948010
$ws.Range("L134").Value = 90948010
$ws.Range("N134").Value = -90958150

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50748
# Row 32
$ws.Range("H32").Value = 4577445
$ws.Range("I32").Value = 6779.6553
$ws.Range("J32").Value = 22250684
$ws.Range("K32").Value = 6779.6553
$ws.Range("L32").Value = 22250684
$ws.Range("M32").Value = -6492.6553
$ws.Range("N32").Value = -22251258
# Row 45
$ws.Range("H45").Value = 2477.8857
$ws.Range("I45").Value = 1442.2778
$ws.Range("J45").Value = 3574.4119
$ws.Range("K45").Value = 1442.2778
$ws.Range("L45").Value = 3574.4119
$ws.Range("M45").Value = -1065.2778
$ws.Range("N45").Value = -4328.4119
# Row 74
$ws.Range("H74").Value = 27274392
$ws.Range("I74").Value = 33335144
$ws.Range("J74").Value = 1002.8
$ws.Range("K74").Value = 33335144
$ws.Range("L74").Value = 1002.8
$ws.Range("M74").Value = -33334270
$ws.Range("N74").Value = -2750.8
# Row 77
$ws.Range("H77").Value = 27274392
$ws.Range("I77").Value = 33335144
$ws.Range("J77").Value = 1002.8
$ws.Range("K77").Value = 166675720
$ws.Range("L77").Value = 5014
$ws.Range("M77").Value = -166671352
$ws.Range("N77").Value = -13750
# Row 100
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164
# Row 122
$ws.Range("H122").Value = 2842.7856
$ws.Range("I122").Value = 1144.4445
$ws.Range("J122").Value = 5899.8
$ws.Range("K122").Value = 3433.3335
$ws.Range("L122").Value = 17699.4
$ws.Range("M122").Value = -983.3335000000002
$ws.Range("N122").Value = -22599.4
# Row 132
$ws.Range("H132").Value = 2454847.5
$ws.Range("I132").Value = 1742.6774
$ws.Range("J132").Value = 6257159.5
$ws.Range("K132").Value = 5228.0322
$ws.Range("L132").Value = 18771478.5
$ws.Range("M132").Value = -2698.0322
$ws.Range("N132").Value = -18776538.5

$ws = $wb.Worksheets.Item("BSM")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 99
$ws.Range("H99").Value = 2712.4
$ws.Range("I99").Value = 2127.5
$ws.Range("K99").Value = 2127.5
$ws.Range("M99").Value = -629.5
# Row 134
$ws.Range("H134").Value = 5027.298
$ws.Range("I134").Value = 2224.8845
$ws.Range("J134").Value = 8496.951999999999
$ws.Range("K134").Value = 6674.6535
$ws.Range("L134").Value = 25490.856
$ws.Range("M134").Value = -4139.6535
$ws.Range("N134").Value = -30560.856

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 10000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -10224
# Row 31
$ws.Range("H31").Value = 7938577.5
$ws.Range("I31").Value = 1526.5294
$ws.Range("J31").Value = 17244086
$ws.Range("K31").Value = 1526.5294
$ws.Range("L31").Value = 17244086
$ws.Range("M31").Value = -1231.5294
$ws.Range("N31").Value = -17244676
# Row 34
$ws.Range("H34").Value = 7938577.5
$ws.Range("I34").Value = 1526.5294
$ws.Range("J34").Value = 17244086
$ws.Range("K34").Value = 1526.5294
$ws.Range("L34").Value = 17244086
$ws.Range("M34").Value = -1324.5294
$ws.Range("N34").Value = -17244490
# Row 74
$ws.Range("H74").Value = 33845
$ws.Range("J74").Value = 33845
$ws.Range("L74").Value = 33845
$ws.Range("N74").Value = -35593
# Row 77
$ws.Range("H77").Value = 33845
$ws.Range("J77").Value = 33845
$ws.Range("L77").Value = 101535
$ws.Range("N77").Value = -110271
# Row 99
$ws.Range("H99").Value = 2599.7083
$ws.Range("I99").Value = 1587.5333
$ws.Range("J99").Value = 4286.6665
$ws.Range("K99").Value = 1587.5333
$ws.Range("L99").Value = 4286.6665
$ws.Range("M99").Value = -89.53330000000005
$ws.Range("N99").Value = -7282.6665
# Row 126
$ws.Range("H126").Value = 2599.7083
$ws.Range("I126").Value = 1587.5333
$ws.Range("J126").Value = 4286.6665
$ws.Range("K126").Value = 4762.5999
$ws.Range("L126").Value = 12859.9995
$ws.Range("M126").Value = -2292.5999
$ws.Range("N126").Value = -17799.9995
# Row 132
$ws.Range("H132").Value = 2191.0513
$ws.Range("I132").Value = 1359.4286
$ws.Range("K132").Value = 4078.2858
$ws.Range("M132").Value = -1548.2858

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 944.33
$ws.Range("I68").Value = 800.4375
$ws.Range("J68").Value = 1077.1538
$ws.Range("K68").Value = 2401.3125
$ws.Range("L68").Value = 3231.4614
$ws.Range("M68").Value = -1590.3125
$ws.Range("N68").Value = -4853.4614
# Row 71
$ws.Range("H71").Value = 944.33
$ws.Range("I71").Value = 800.4375
$ws.Range("J71").Value = 1077.1538
$ws.Range("K71").Value = 7203.9375
$ws.Range("L71").Value = 9694.3842
$ws.Range("M71").Value = -3147.9375
$ws.Range("N71").Value = -17806.3842

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2581.9583
$ws.Range("I122").Value = 1984.579
$ws.Range("J122").Value = 4852
$ws.Range("K122").Value = 5953.737
$ws.Range("L122").Value = 14556
$ws.Range("M122").Value = -3503.737
$ws.Range("N122").Value = -19456

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 34523144
$ws.Range("I132").Value = 66739600
$ws.Range("J132").Value = 5512.143
$ws.Range("K132").Value = 200218800
$ws.Range("L132").Value = 16536.429
$ws.Range("M132").Value = -200216270
$ws.Range("N132").Value = -21596.429

$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 9000
$ws.Range("I76").Value = 9000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -8685
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 9000
$ws.Range("I79").Value = 9000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -7908
$ws.Range("N79").ClearContents()
# Row 132
$ws.Range("H132").Value = 2285.558
$ws.Range("I132").Value = 1748.4762
$ws.Range("J132").Value = 2798.2273
$ws.Range("K132").Value = 5245.4286
$ws.Range("L132").Value = 8394.6819
$ws.Range("M132").Value = -2715.4286
$ws.Range("N132").Value = -13454.6819
